$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '88.981.85'
$ws.Range("E2").Value = '  -3.04%  '

$ws.Range("D3").Value = '3.090.60'
$ws.Range("E3").Value = '  -6.43%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.07'
$ws.Range("E5").Value = '  -1.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '626.19'
$ws.Range("E6").Value = '  -0.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.386'
$ws.Range("E7").Value = '  -8.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.711'
$ws.Range("E8").Value = '  -0.35%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").Value = '3.086.62'
$ws.Range("E10").Value = '  -6.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.548'
$ws.Range("E11").Value = '  -7.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.178'
$ws.Range("E12").Value = '  -0.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000247'
$ws.Range("E13").Value = '  -8.15%  '

$ws.Range("D14").Value = '88.660.51'
$ws.Range("E14").Value = '  -2.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.21'
$ws.Range("E15").Value = '  -3.86%  '

$ws.Range("D16").Value = '3.668.83'
$ws.Range("E16").Value = '  -5.96%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.68'
$ws.Range("E17").Value = '  -8.22%  '

$ws.Range("D18").Value = '3.086.88'
$ws.Range("E18").Value = '  -6.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.30'
$ws.Range("E19").Value = '  +0.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000205'
$ws.Range("E20").Value = '  +8.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.96'
$ws.Range("E21").Value = '  -8.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '421.44'
$ws.Range("E22").Value = '  -3.85%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.20'
$ws.Range("E23").Value = '  -9.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.81'
$ws.Range("E24").Value = '  -10.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.21'
$ws.Range("E25").Value = '  -4.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.16'
$ws.Range("E26").Value = '  -8.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '78.31'
$ws.Range("E27").Value = '  +2.16%  '

$ws.Range("D28").Value = '3.292.15'
$ws.Range("E28").Value = '  -5.02%  '

$ws.Range("E29").Value = '  +0.20%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("E31").Value = '  -11.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.08'
$ws.Range("E32").Value = '  -7.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.82'
$ws.Range("E33").Value = '  +1.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '499.67'
$ws.Range("E34").Value = '  -11.17%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.71'
$ws.Range("E35").Value = '  -9.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.81'
$ws.Range("E36").Value = '  -6.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.24'
$ws.Range("E37").Value = '  -7.97%  '

$ws.Range("B38").Value = 'WhiteBITCoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.19'
$ws.Range("E38").Value = '  -0.77%  '

$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.59'
$ws.Range("E39").Value = '  -5.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.123'
$ws.Range("E41").Value = '  -8.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.84'
$ws.Range("E43").Value = '  -8.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.358'
$ws.Range("E44").Value = '  -9.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '144.24'
$ws.Range("E45").Value = '  -3.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.31'
$ws.Range("E46").Value = '  -1.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '164.10'
$ws.Range("E47").Value = '  -10.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.123'
$ws.Range("E48").Value = '  -5.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.709'
$ws.Range("E49").Value = '  -3.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.83'
$ws.Range("E50").Value = '  -9.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.16'
$ws.Range("E51").Value = '  -10.26%  '
